# Seabird transects workbook update
# - Insert a new "Transect Sequence - Side of Motu" column (G), shifting the
#   old "Transect name (alt)"/"notes" columns one slot to the right.
# - Insert a new "notes_side_of_motu" column (K) at the end.
# - Rename the (now shifted) notes column header to "notes_100-200m_transects".
# - Populate the new columns with their data.
# - Resize rows whose wrapped text now needs more vertical space.
# - Update window/selection state to reflect where the author left off (G5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stage a "bold, not wrapped" header style on a scratch cell far outside
#     the table so we can clone it onto the new K1 header without creating a
#     brand-new (duplicate) style entry in the workbook. ---
$ws.Range("Z1").Value = "scratch"
$ws.Range("Z1").Font.Bold = $true

# --- Insert the two new columns ---
# New column G: "Transect Sequence - Side of Motu"
$ws.Columns("G").Insert()
# New column K: "notes_side_of_motu" (after the shifted notes column J)
$ws.Columns("K").Insert()

# --- Formatting for the new columns ---
# Column G should look like the other data columns (C:F / now shifted to match F)
$ws.Range("F1:F7").Copy()
$ws.Range("G1:G7").PasteSpecial(-4122)   # xlPasteFormats

# Column K data rows (2:7) should look like column J (old plain notes style)
$ws.Range("J2:J7").Copy()
$ws.Range("K2:K7").PasteSpecial(-4122)   # xlPasteFormats

# Column K header (K1) should be bold / not wrapped - clone from scratch cell
# (it shifted two columns to the right, from Z1 to AB1, after the two inserts)
$ws.Range("AB1").Copy()
$ws.Range("K1").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = 0

# Remove the scratch column entirely now that its format has been cloned
$ws.Columns("AB").Delete()

# --- Header text updates ---
$ws.Range("G1").Value = "Transect Sequence - Side of Motu"
$ws.Range("J1").Value = "notes_100-200m_transects"
$ws.Range("K1").Value = "notes_side_of_motu"

# --- New "Transect Sequence - Side of Motu" values ---
$ws.Range("G2").Value = "0945-0946-0947-0948"
$ws.Range("G3").Value = "0948-0949-0950-0945"
$ws.Range("G4").Value = "0692 - 0693 - 0694 - 0695 - 0696 - 0697 - 0698 - 0699 - 0700 - 0701 - 0702 - 0606"
$ws.Range("G5").Value = "0631 - 0632 - 0633 - 0634 - 0635 - 0636 - 0637 - 0638 - 0644 - 0645 - 0646 - 0647 - 0648 - 0686 - 0687 - 0688 - 0689 - 0690 - 0691  - 0692"
$ws.Range("G6").Value = "1143-1144-1145-1146-1147-1148-1149-1150-1151-1152"
$ws.Range("G7").Value = "0837-0838-0839-1137-1138-1139-1140-1141-1142-1143"

# --- New "notes_side_of_motu" values ---
$ws.Range("K2").Value = "obvious split to e/w corners of motu, just added 1 extra transect to 200 m"
$ws.Range("K3").Value = "obvious split to e/w corners of motu, just added 1 extra transect to 200 m"
$ws.Range("K4").Value = "only included until eastern most corner (so excluded northern stretch which is cut-off from transects). Obvious southern corner. "
$ws.Range("K5").Value = "only included until northwestern corner (so excluded norther stretch which is cut-off from transects). Obvious southern corner"
$ws.Range("K6").Value = "obvious split to n/s corners of motu"
$ws.Range("K7").Value = "obvious split to n/s corners of motu"

# --- Row heights now need to grow to fit the newly-added wrapped text ---
$ws.Rows(2).RowHeight = 68
$ws.Rows(3).RowHeight = 68
$ws.Rows(4).RowHeight = 136
$ws.Rows(5).RowHeight = 119

# --- Selection moves to G5 (where the author was last working) ---
$ws.Range("G5").Select()
